$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '69.456.03'
$ws.Range('E2').Value = '  +0.20%  '
$ws.Range('D3').Value = '3.695.02'
$ws.Range('E3').Value = '  +0.33%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '0.999'
$ws.Range('D4').Style = 'Normal'
$ws.Range('E4').Value = '  -0.08%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '681.25'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -0.25%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '159.72'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -0.26%  '
$ws.Range('E7').Value = '  -0.09%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.495'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  +0.21%  '
$ws.Range('E9').Value = '  +0.47%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '7.15'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -1.19%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.441'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +1.35%  '
$ws.Range('E12').Value = '  -0.67%  '
$ws.Range('D13').Value = '4.315.20'
$ws.Range('E13').Value = '  +0.26%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '32.36'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -0.65%  '
$ws.Range('D15').Value = '3.686.86'
$ws.Range('E15').Value = '  +0.11%  '
$ws.Range('D16').Value = '69.408.14'
$ws.Range('E16').Value = '  +0.05%  '
$ws.Range('E17').Value = '  +3.12%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '16.02'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  +1.18%  '
$ws.Range('E19').Value = '  +0.56%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '469.51'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -0.24%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '9.92'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +0.40%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '0.652'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +0.10%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '80.21'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +0.77%  '
$ws.Range('D24').Value = '3.839.57'
$ws.Range('E24').Value = '  +0.28%  '
$ws.Range('E25').Value = '  -0.06%  '
$ws.Range('E26').Value = '  -2.49%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '10.94'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -0.85%  '
$ws.Range('E28').Value = '  +0.38%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '2.71'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  +0.21%  '
$ws.Range('E30').Value = '  -1.04%  '
$ws.Range('B31').Value = 'ImmutableX'
$ws.Range('C31').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '1.99'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -1.49%  '
$ws.Range('B32').Value = 'Binance-PegBSC-USD'
$ws.Range('C32').Value = 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '1.00'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +0.07%  '
$ws.Range('B33').Value = 'NEARProtocol'
$ws.Range('C33').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '6.55'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -2.29%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '27.01'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +1.15%  '
$ws.Range('D35').Value = '3.683.55'
$ws.Range('E35').Value = '  +0.85%  '
$ws.Range('E36').Value = '  -1.32%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '8.34'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +2.07%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '6.28'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +2.27%  '
$ws.Range('B40').Value = 'Stacks'
$ws.Range('C40').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '2.24'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -1.10%  '
$ws.Range('B41').Value = 'FirstDigitalUSD'
$ws.Range('C41').Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.999'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -0.16%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.0905'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +0.10%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '169.42'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +2.70%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.943'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -0.14%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '47.05'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -1.85%  '
$ws.Range('B46').Value = 'dogwifhat'
$ws.Range('C46').Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '2.71'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -0.42%  '
$ws.Range('B47').Value = 'FLOKI'
$ws.Range('C47').Value = 'https://coinranking.com/coin/fmHk13Rqw+floki-floki'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.000279'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +2.13%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '28.06'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -3.32%  '
$ws.Range('B49').Value = 'SuiNetwork'
$ws.Range('C49').Value = 'https://coinranking.com/coin/3xJluUMvp+suinetwork-sui'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '1.10'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +0.64%  '
$ws.Range('B50').Value = 'ONDO'
$ws.Range('C50').Value = 'https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '1.29'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -1.37%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '7.83'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -0.59%  '

Write-Host "Applied cryptos update"
